# Plant Promoter Comparison spreadsheet update
# - Rename column K header from "Name and location of Fastq Reads" to
#   "Name and location of Original Reads"
# - Add two new data rows (6 and 7) describing the Arabidopsis
#   EST / cDNA original-read sources
# - Update the window/sheet view position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# --- Header update -------------------------------------------------------
$ws.Range("K1").Value = "Name and location of Original Reads"

# --- New row 6 : Arabidopsis thaliana / EST ------------------------------
$ws.Range("A6").Value = "Arabidopsis thaliana"
$ws.Range("A6").Font.Italic = $true
$ws.Range("B6").Value = "Mutliple/NA"
$ws.Range("C6").Value = "EST"
$ws.Range("D6").Value = 1816638
$ws.Range("E6").Value = "N/A"
$ws.Range("F6").Value = 1
$ws.Range("J6").Value = "/scratch/rtraborn/TSRchitect_plant_results/Arabidopsis/EST_cDNA/TH_EST_sequences_20101108.gsq "
$ws.Range("K6").Value = "/scratch/rtraborn/TSRchitect_plant_results/Arabidopsis/EST_cDNA/ATH_EST_sequences_20101108.fas"

# --- New row 7 : Arabidopsis thaliana / cDNA -----------------------------
$ws.Range("A7").Value = "Arabidopsis thaliana"
$ws.Range("A7").Font.Italic = $true
$ws.Range("B7").Value = "Mutliple/NA"
$ws.Range("C7").Value = "cDNA"
$ws.Range("D7").Value = 78096
$ws.Range("E7").Value = "N/A"
$ws.Range("F7").Value = 1
$ws.Range("J7").Value = "/scratch/rtraborn/TSRchitect_plant_results/Arabidopsis/EST_cDNA/ATH_cDNA_sequences_20101108.gsq "
$ws.Range("K7").Value = "/scratch/rtraborn/TSRchitect_plant_results/Arabidopsis/EST_cDNA/ATH_cDNA_sequences_20101108.fas"

# --- Window / selection state ---------------------------------------------
$ws.Activate()
$win = $wb.Windows.Item(1)
$win.Left = 51700
$win.Top = 1400
$ws.Range("F13").Select()
